$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8769673109054565
$ws.Range("B1").Value = 1.308611154556274
$ws.Range("D1").Value = 1.702446818351746
$ws.Range("E1").Value = 1.112001776695251
